$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9146688338073954
$ws.Range("D3").Value = 0.08533116619260463
$ws.Range("D4").Value = 2251
$ws.Range("D5").Value = 210
$ws.Range("D6").Value = 0.977491961414791
$ws.Range("D7").Value = 0.022508038585209
$ws.Range("D8").Value = 304
$ws.Range("D9").Value = 7
